$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 5.305917
$ws.Range("H2").Value = 15.917751
$ws.Range("I2").Value = 0.4336744870332215
$ws.Range("J2").Value = 0.4336744870332215
$ws.Range("M2").Value = 0.08849299999999999
$ws.Range("Q2").Value = 0.4695365130809999
$ws.Range("R2").Value = 4.225828617728999
$ws.Range("S2").Value = 0.392575728822024
$ws.Range("T2").Value = 0.3925757288220239

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 5.305917
$ws.Range("H3").Value = 15.917751
$ws.Range("I3").Value = 0.4336744870332215
$ws.Range("J3").Value = 0.4336744870332215
$ws.Range("O3").Value = 0.09476867890558938
$ws.Range("P3").Value = 0.09476867890558936
$ws.Range("Q3").Value = 0.049155783727
$ws.Range("R3").Value = 0.442402053543
$ws.Range("S3").Value = 0.04109875821119756
$ws.Range("T3").Value = 0.04109875821119754

# Row 4
$ws.Range("I4").Value = 0.2066699405724794
$ws.Range("J4").Value = 0.2066699405724794
$ws.Range("M4").Value = 0.08849299999999999
$ws.Range("Q4").Value = 0.2237601845543333
$ws.Range("S4").Value = 0.1870841033349289
$ws.Range("T4").Value = 0.1870841033349289

# Row 5
$ws.Range("I5").Value = 0.2066699405724794
$ws.Range("J5").Value = 0.2066699405724794
$ws.Range("O5").Value = 0.09476867890558938
$ws.Range("P5").Value = 0.09476867890558936
$ws.Range("S5").Value = 0.01958583723755054
$ws.Range("T5").Value = 0.01958583723755054

# Row 6
$ws.Range("G6").Value = 4.400310999999999
$ws.Range("H6").Value = 13.200933
$ws.Range("I6").Value = 0.359655572394299
$ws.Range("J6").Value = 0.359655572394299
$ws.Range("M6").Value = 0.08849299999999999
$ws.Range("Q6").Value = 0.3893967213229999
$ws.Range("R6").Value = 3.504570491906999
$ws.Range("S6").Value = 0.3255714889374577
$ws.Range("T6").Value = 0.3255714889374577

# Row 7
$ws.Range("G7").Value = 4.400310999999999
$ws.Range("H7").Value = 13.200933
$ws.Range("I7").Value = 0.359655572394299
$ws.Range("J7").Value = 0.359655572394299
$ws.Range("O7").Value = 0.09476867890558938
$ws.Range("P7").Value = 0.09476867890558936
$ws.Range("Q7").Value = 0.04076594787433333
$ws.Range("R7").Value = 0.366893530869
$ws.Range("S7").Value = 0.03408408345684128
$ws.Range("T7").Value = 0.03408408345684127
